$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing cell CONTENT only (keeps formatting/styles.xml untouched)
# and resets the shared-string table so we can rebuild it in the desired order.
$ws.Range("A1:M19").ClearContents()

# --- Row 1: column index headers (B1:M1 = 0..11) ---
$ws.Cells.Item(1,2).Value2 = 0
$ws.Cells.Item(1,3).Value2 = 1
$ws.Cells.Item(1,4).Value2 = 2
$ws.Cells.Item(1,5).Value2 = 3
$ws.Cells.Item(1,6).Value2 = 4
$ws.Cells.Item(1,7).Value2 = 5
$ws.Cells.Item(1,8).Value2 = 6
$ws.Cells.Item(1,9).Value2 = 7
$ws.Cells.Item(1,10).Value2 = 8
$ws.Cells.Item(1,11).Value2 = 9
$ws.Cells.Item(1,12).Value2 = 10
$ws.Cells.Item(1,13).Value2 = 11

# --- Column A (index numbers) and Column B (labels), written row by row,
# --- BEFORE the row-2 header labels, so the label strings occupy the first
# --- shared-string slots and the bracket/header labels are appended after them. ---
$ws.Cells.Item(2,1).Value2 = 0
$ws.Cells.Item(2,2).Value2 = "HKL"
$ws.Cells.Item(3,1).Value2 = 1
$ws.Cells.Item(3,2).Value2 = "ND Single"
$ws.Cells.Item(4,1).Value2 = 2
$ws.Cells.Item(4,2).Value2 = "RD Single"
$ws.Cells.Item(5,1).Value2 = 3
$ws.Cells.Item(5,2).Value2 = "TD Single"
$ws.Cells.Item(6,1).Value2 = 4
$ws.Cells.Item(6,2).Value2 = "Morris"
$ws.Cells.Item(7,1).Value2 = 5
$ws.Cells.Item(7,2).Value2 = "Ring Perpendicular to ND"
$ws.Cells.Item(8,1).Value2 = 6
$ws.Cells.Item(8,2).Value2 = "Ring Perpendicular to RD"
$ws.Cells.Item(9,1).Value2 = 7
$ws.Cells.Item(9,2).Value2 = "Ring Perpendicular to TD"
$ws.Cells.Item(10,1).Value2 = 8
$ws.Cells.Item(10,2).Value2 = "Gaussian-Quadrature"
$ws.Cells.Item(11,1).Value2 = 9
$ws.Cells.Item(11,2).Value2 = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(12,1).Value2 = 10
$ws.Cells.Item(12,2).Value2 = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(13,1).Value2 = 11
$ws.Cells.Item(13,2).Value2 = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(14,1).Value2 = 12
$ws.Cells.Item(14,2).Value2 = "NoRotation-tilt60deg"
$ws.Cells.Item(15,1).Value2 = 13
$ws.Cells.Item(15,2).Value2 = "Rotation-NoTilt"
$ws.Cells.Item(16,1).Value2 = 14
$ws.Cells.Item(16,2).Value2 = "Rotation-60detTilt"
$ws.Cells.Item(17,1).Value2 = 15
$ws.Cells.Item(17,2).Value2 = "HexGrid-90degTilt5degRes"
$ws.Cells.Item(18,1).Value2 = 16
$ws.Cells.Item(18,2).Value2 = "HexGrid-90degTilt22p5degRes"
$ws.Cells.Item(19,1).Value2 = 17
$ws.Cells.Item(19,2).Value2 = "HexGrid-60degTilt5degRes"

# --- Row 2 header labels C2:M2 (bracket / pairs notation), introduced AFTER all B-column labels ---
$ws.Cells.Item(2,3).Value2 = "[1, 1, 0]"
$ws.Cells.Item(2,4).Value2 = "[2, 0, 0]"
$ws.Cells.Item(2,5).Value2 = "[2, 1, 1]"
$ws.Cells.Item(2,6).Value2 = "[2, 2, 0]"
$ws.Cells.Item(2,7).Value2 = "[3, 1, 0]"
$ws.Cells.Item(2,8).Value2 = "[2, 2, 2]"
$ws.Cells.Item(2,9).Value2 = "[3, 2, 1]"
$ws.Cells.Item(2,10).Value2 = "[4, 0, 0]"
$ws.Cells.Item(2,11).Value2 = "2Pairs"
$ws.Cells.Item(2,12).Value2 = "4Pairs"
$ws.Cells.Item(2,13).Value2 = "MaxUnique"

# --- Numeric data for rows 3-19, columns C:M ---
$ws.Cells.Item(3,3).Value2 = 0.69
$ws.Cells.Item(3,4).Value2 = 0.1
$ws.Cells.Item(3,5).Value2 = 1.51
$ws.Cells.Item(3,6).Value2 = 0.69
$ws.Cells.Item(3,7).Value2 = 0.26
$ws.Cells.Item(3,8).Value2 = 2.94
$ws.Cells.Item(3,9).Value2 = 1.36
$ws.Cells.Item(3,10).Value2 = 0.1
$ws.Cells.Item(3,11).Value2 = 0.805
$ws.Cells.Item(3,12).Value2 = 0.7474999999999999
$ws.Cells.Item(3,13).Value2 = 1.143333333333333
$ws.Cells.Item(4,3).Value2 = 1.41
$ws.Cells.Item(4,4).Value2 = 0.35
$ws.Cells.Item(4,5).Value2 = 1.06
$ws.Cells.Item(4,6).Value2 = 1.41
$ws.Cells.Item(4,7).Value2 = 0.66
$ws.Cells.Item(4,8).Value2 = 1.11
$ws.Cells.Item(4,9).Value2 = 1.16
$ws.Cells.Item(4,10).Value2 = 0.35
$ws.Cells.Item(4,11).Value2 = 0.7050000000000001
$ws.Cells.Item(4,12).Value2 = 1.0575
$ws.Cells.Item(4,13).Value2 = 0.9583333333333334
$ws.Cells.Item(5,3).Value2 = 1.41
$ws.Cells.Item(5,4).Value2 = 0.35
$ws.Cells.Item(5,5).Value2 = 1.06
$ws.Cells.Item(5,6).Value2 = 1.41
$ws.Cells.Item(5,7).Value2 = 0.67
$ws.Cells.Item(5,8).Value2 = 1.11
$ws.Cells.Item(5,9).Value2 = 1.16
$ws.Cells.Item(5,10).Value2 = 0.35
$ws.Cells.Item(5,11).Value2 = 0.7050000000000001
$ws.Cells.Item(5,12).Value2 = 1.0575
$ws.Cells.Item(5,13).Value2 = 0.96
$ws.Cells.Item(6,3).Value2 = 0.65
$ws.Cells.Item(6,4).Value2 = 1.58
$ws.Cells.Item(6,5).Value2 = 0.98
$ws.Cells.Item(6,6).Value2 = 0.65
$ws.Cells.Item(6,7).Value2 = 1.23
$ws.Cells.Item(6,8).Value2 = 0.96
$ws.Cells.Item(6,9).Value2 = 0.89
$ws.Cells.Item(6,10).Value2 = 1.58
$ws.Cells.Item(6,11).Value2 = 1.28
$ws.Cells.Item(6,12).Value2 = 0.965
$ws.Cells.Item(6,13).Value2 = 1.048333333333333
$ws.Cells.Item(7,3).Value2 = 1.39986301369863
$ws.Cells.Item(7,4).Value2 = 0.35
$ws.Cells.Item(7,5).Value2 = 1.051780821917808
$ws.Cells.Item(7,6).Value2 = 1.39986301369863
$ws.Cells.Item(7,7).Value2 = 0.6647945205479452
$ws.Cells.Item(7,8).Value2 = 1.11
$ws.Cells.Item(7,9).Value2 = 1.16
$ws.Cells.Item(7,10).Value2 = 0.35
$ws.Cells.Item(7,11).Value2 = 0.700890410958904
$ws.Cells.Item(7,12).Value2 = 1.050376712328767
$ws.Cells.Item(7,13).Value2 = 0.9560730593607306
$ws.Cells.Item(8,3).Value2 = 0.9810526315789474
$ws.Cells.Item(8,4).Value2 = 0.8121052631578948
$ws.Cells.Item(8,5).Value2 = 1.081052631578947
$ws.Cells.Item(8,6).Value2 = 0.9810526315789474
$ws.Cells.Item(8,7).Value2 = 0.8642105263157894
$ws.Cells.Item(8,8).Value2 = 1.268421052631579
$ws.Cells.Item(8,9).Value2 = 1.055789473684211
$ws.Cells.Item(8,10).Value2 = 0.8121052631578948
$ws.Cells.Item(8,11).Value2 = 0.946578947368421
$ws.Cells.Item(8,12).Value2 = 0.9638157894736843
$ws.Cells.Item(8,13).Value2 = 1.010438596491228
$ws.Cells.Item(9,3).Value2 = 0.9810526315789474
$ws.Cells.Item(9,4).Value2 = 0.8121052631578948
$ws.Cells.Item(9,5).Value2 = 1.081052631578947
$ws.Cells.Item(9,6).Value2 = 0.9810526315789474
$ws.Cells.Item(9,7).Value2 = 0.8642105263157894
$ws.Cells.Item(9,8).Value2 = 1.268421052631579
$ws.Cells.Item(9,9).Value2 = 1.055789473684211
$ws.Cells.Item(9,10).Value2 = 0.8121052631578948
$ws.Cells.Item(9,11).Value2 = 0.946578947368421
$ws.Cells.Item(9,12).Value2 = 0.9638157894736843
$ws.Cells.Item(9,13).Value2 = 1.010438596491228
$ws.Cells.Item(10,3).Value2 = 1.186928358257481
$ws.Cells.Item(10,4).Value2 = 0.6025458675910295
$ws.Cells.Item(10,5).Value2 = 1.05201901139688
$ws.Cells.Item(10,6).Value2 = 1.186928358257481
$ws.Cells.Item(10,7).Value2 = 0.7990077946846026
$ws.Cells.Item(10,8).Value2 = 1.1312998136103
$ws.Cells.Item(10,9).Value2 = 1.09943573977727
$ws.Cells.Item(10,10).Value2 = 0.6025458675910295
$ws.Cells.Item(10,11).Value2 = 0.8272824394939546
$ws.Cells.Item(10,12).Value2 = 1.007105398875717
$ws.Cells.Item(10,13).Value2 = 0.9785394308862604
$ws.Cells.Item(11,3).Value2 = 0.9933438884412549
$ws.Cells.Item(11,4).Value2 = 0.8805574828700983
$ws.Cells.Item(11,5).Value2 = 1.047344888658173
$ws.Cells.Item(11,6).Value2 = 0.9933438884412549
$ws.Cells.Item(11,7).Value2 = 0.9168840163964505
$ws.Cells.Item(11,8).Value2 = 1.143854903952787
$ws.Cells.Item(11,9).Value2 = 1.031188684351605
$ws.Cells.Item(11,10).Value2 = 0.8805574828700983
$ws.Cells.Item(11,11).Value2 = 0.9639511857641354
$ws.Cells.Item(11,12).Value2 = 0.9786475371026953
$ws.Cells.Item(11,13).Value2 = 1.002195644111728
$ws.Cells.Item(12,3).Value2 = 0.9925108790721635
$ws.Cells.Item(12,4).Value2 = 0.8815462107277154
$ws.Cells.Item(12,5).Value2 = 1.047339437449253
$ws.Cells.Item(12,6).Value2 = 0.9925108790721635
$ws.Cells.Item(12,7).Value2 = 0.9173862341497929
$ws.Cells.Item(12,8).Value2 = 1.143508101440162
$ws.Cells.Item(12,9).Value2 = 1.030761335621283
$ws.Cells.Item(12,10).Value2 = 0.8815462107277154
$ws.Cells.Item(12,11).Value2 = 0.9644428240884841
$ws.Cells.Item(12,12).Value2 = 0.9784768515803239
$ws.Cells.Item(12,13).Value2 = 1.002175366410062
$ws.Cells.Item(13,3).Value2 = 0.9934469166709441
$ws.Cells.Item(13,4).Value2 = 0.8803512120821483
$ws.Cells.Item(13,5).Value2 = 1.047368256941511
$ws.Cells.Item(13,6).Value2 = 0.9934469166709441
$ws.Cells.Item(13,7).Value2 = 0.9168159963610443
$ws.Cells.Item(13,8).Value2 = 1.143335735184324
$ws.Cells.Item(13,9).Value2 = 1.031013480838639
$ws.Cells.Item(13,10).Value2 = 0.8803512120821483
$ws.Cells.Item(13,11).Value2 = 0.9638597345118298
$ws.Cells.Item(13,12).Value2 = 0.978653325591387
$ws.Cells.Item(13,13).Value2 = 1.002055266346435
$ws.Cells.Item(14,3).Value2 = 0.8264320000000005
$ws.Cells.Item(14,4).Value2 = 0.4700000000000004
$ws.Cells.Item(14,5).Value2 = 1.295755999999998
$ws.Cells.Item(14,6).Value2 = 0.8264320000000005
$ws.Cells.Item(14,7).Value2 = 0.5663559999999996
$ws.Cells.Item(14,8).Value2 = 2.107623999999998
$ws.Cells.Item(14,9).Value2 = 1.203604
$ws.Cells.Item(14,10).Value2 = 0.4700000000000004
$ws.Cells.Item(14,11).Value2 = 0.8828779999999993
$ws.Cells.Item(14,12).Value2 = 0.8546549999999999
$ws.Cells.Item(14,13).Value2 = 1.078295333333333
$ws.Cells.Item(15,3).Value2 = 0.69
$ws.Cells.Item(15,4).Value2 = 0.1
$ws.Cells.Item(15,5).Value2 = 1.51
$ws.Cells.Item(15,6).Value2 = 0.69
$ws.Cells.Item(15,7).Value2 = 0.26
$ws.Cells.Item(15,8).Value2 = 2.94
$ws.Cells.Item(15,9).Value2 = 1.36
$ws.Cells.Item(15,10).Value2 = 0.1
$ws.Cells.Item(15,11).Value2 = 0.805
$ws.Cells.Item(15,12).Value2 = 0.7474999999999999
$ws.Cells.Item(15,13).Value2 = 1.143333333333333
$ws.Cells.Item(16,3).Value2 = 0.8263879137280029
$ws.Cells.Item(16,4).Value2 = 0.4695855489024012
$ws.Cells.Item(16,5).Value2 = 1.295642868940795
$ws.Cells.Item(16,6).Value2 = 0.8263879137280029
$ws.Cells.Item(16,7).Value2 = 0.5664575488000008
$ws.Cells.Item(16,8).Value2 = 2.107015674675212
$ws.Cells.Item(16,9).Value2 = 1.203729755135997
$ws.Cells.Item(16,10).Value2 = 0.4695855489024012
$ws.Cells.Item(16,11).Value2 = 0.8826142089215983
$ws.Cells.Item(16,12).Value2 = 0.8545010613248005
$ws.Cells.Item(16,13).Value2 = 1.078136551697068
$ws.Cells.Item(17,3).Value2 = 0.9957089400709754
$ws.Cells.Item(17,4).Value2 = 0.9980709158298056
$ws.Cells.Item(17,5).Value2 = 1.000374582478758
$ws.Cells.Item(17,6).Value2 = 0.9957089400709754
$ws.Cells.Item(17,7).Value2 = 0.9958305039463997
$ws.Cells.Item(17,8).Value2 = 1.00115640246909
$ws.Cells.Item(17,9).Value2 = 0.9987901497843227
$ws.Cells.Item(17,10).Value2 = 0.9980709158298056
$ws.Cells.Item(17,11).Value2 = 0.9992227491542817
$ws.Cells.Item(17,12).Value2 = 0.9974658446126285
$ws.Cells.Item(17,13).Value2 = 0.9983219157632252
$ws.Cells.Item(18,3).Value2 = 0.967419753829643
$ws.Cells.Item(18,4).Value2 = 1.055795738000985
$ws.Cells.Item(18,5).Value2 = 0.9929893569424617
$ws.Cells.Item(18,6).Value2 = 0.967419753829643
$ws.Cells.Item(18,7).Value2 = 1.024547616540342
$ws.Cells.Item(18,8).Value2 = 0.989985832573801
$ws.Cells.Item(18,9).Value2 = 0.9866461243156752
$ws.Cells.Item(18,10).Value2 = 1.055795738000985
$ws.Cells.Item(18,11).Value2 = 1.024392547471723
$ws.Cells.Item(18,12).Value2 = 0.9959061506506832
$ws.Cells.Item(18,13).Value2 = 1.002897403700485
$ws.Cells.Item(19,3).Value2 = 0.9765100825094638
$ws.Cells.Item(19,4).Value2 = 1.160252217770952
$ws.Cells.Item(19,5).Value2 = 0.9536619033252548
$ws.Cells.Item(19,6).Value2 = 0.9765100825094638
$ws.Cells.Item(19,7).Value2 = 1.087040054211242
$ws.Cells.Item(19,8).Value2 = 0.8699926445119092
$ws.Cells.Item(19,9).Value2 = 0.9544050439854842
$ws.Cells.Item(19,10).Value2 = 1.160252217770952
$ws.Cells.Item(19,11).Value2 = 1.056957060548104
$ws.Cells.Item(19,12).Value2 = 1.016733571528784
$ws.Cells.Item(19,13).Value2 = 1.000310324385718

# --- Re-apply the existing "A column" cell format to the newly-added cells A17:A19 ---
# (A1:A16 already retained their original style since only their content was cleared)
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
